$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 10750.2153668132
$ws.Range("F2").Value = -6.27802003238821

# Row 3
$ws.Range("C3").Value = 10233.679553714
$ws.Range("F3").Value = 323.998017416478

# Row 4
$ws.Range("C4").Value = 10612.4695438593
$ws.Range("F4").Value = 341.89575723583

# Row 5
$ws.Range("C5").Value = 10953.7872131484
$ws.Range("F5").Value = 362.850522797411

# Row 6
$ws.Range("C6").Value = 10401.5132426464
$ws.Range("F6").Value = 341.77479902718

# Row 7
$ws.Range("B7").Value = 4485.60099223776
$ws.Range("C7").Value = 7178.80116541052
$ws.Range("E7").Value = 6068.59653844236
$ws.Range("F7").Value = 176.833237660537
